# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# For D-column cells whose new text looks numeric (e.g. "1.00", "302.27",
# "2.302.19"), force the Text number format first so Excel's Value setter
# stores the literal string instead of silently re-parsing it as a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.034.03"
$ws.Range("E2").Value = "  +2.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.302.19"
$ws.Range("E3").Value = "  +2.09%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.27"
$ws.Range("E5").Value = "  +1.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.19"
$ws.Range("E6").Value = "  +6.31%  "

$ws.Range("E7").Value = "  +2.07%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +3.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.34"
$ws.Range("E10").Value = "  +4.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0798"
$ws.Range("E11").Value = "  +1.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.13"
$ws.Range("E12").Value = "  +3.27%  "

$ws.Range("E13").Value = "  +4.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.10"
$ws.Range("E14").Value = "  +18.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.79"
$ws.Range("E15").Value = "  +2.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.659.12"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.318.34"
$ws.Range("E17").Value = "  +3.25%  "

$ws.Range("E18").Value = "  +5.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.947.41"
$ws.Range("E19").Value = "  +2.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.39"
$ws.Range("E20").Value = "  +9.64%  "

$ws.Range("D21").Value = "0.0₃0904"
$ws.Range("E21").Value = "  +1.75%  "

$ws.Range("E22").Value = "  +2.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.91"
$ws.Range("E23").Value = "  +2.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.46"
$ws.Range("E24").Value = "  +1.49%  "

$ws.Range("E25").Value = "  +15.27%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.46"
$ws.Range("E27").Value = "  +0.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.74"
$ws.Range("E28").Value = "  +4.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "168.76"
$ws.Range("E29").Value = "  +0.65%  "

$ws.Range("E30").Value = "  -3.82%  "

$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.16"
$ws.Range("E31").Value = "  +1.66%  "

$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.70"
$ws.Range("E32").Value = "  +0.77%  "

$ws.Range("E33").Value = "  -0.11%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.04"
$ws.Range("E34").Value = "  +2.64%  "

$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").Value = "  +3.88%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("E36").Value = "  +1.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.94"
$ws.Range("E37").Value = "  +6.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0700"
$ws.Range("E38").Value = "  +1.53%  "

$ws.Range("E39").Value = "  +3.85%  "

$ws.Range("E40").Value = "  +5.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.80"
$ws.Range("E41").Value = "  +1.11%  "

$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("E43").Value = "  -2.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.995.99"
$ws.Range("E44").Value = "  +2.55%  "

$ws.Range("E45").Value = "  +3.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.01"
$ws.Range("E46").Value = "  +5.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.56"
$ws.Range("E47").Value = "  +2.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.86"
$ws.Range("E48").Value = "  +3.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.50"
$ws.Range("E49").Value = "  +7.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.530.32"
$ws.Range("E50").Value = "  +2.00%  "

$ws.Range("E51").Value = "  +3.60%  "
